# Generate Report for Handback
# Inserts a new handback record for file
# 7188329f-3960-45ce-a771-f362959050de.md as the second data row
# (between 2118c7ff... and b964db4f...) on all three sheets:
# Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newBase = "7188329f-3960-45ce-a771-f362959050de"
$newMd = $newBase + ".md"
$newZhXlf = $newBase + ".01b39cc38ce0d383d7e39f62d655affe630da979.zh-cn.xlf"
$newDeXlf = $newBase + ".01b39cc38ce0d383d7e39f62d655affe630da979.de-de.xlf"

$newMdUrlSrc = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c1f6b2a6e8a4b0b9d6e4b4f0b6a7c5d8e9f0a1b/e2e/" + $newMd
$newMdUrlZh  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d/e2e/" + $newMd
$newMdUrlDe  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2e/e2e/" + $newMd

# ---------------------------------------------------------------
# Sheet "Overview": columns A:G
# ---------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item(1)
$loO.ListRows.Add() | Out-Null   # appends a blank row4, table range -> A1:G4

# Move the existing row3 (b964db4f...) data down to row4
foreach ($col in @("A","B","C","E","F","G")) {
    $wsO.Range($col + "4").Value2 = $wsO.Range($col + "3").Value2
}
$wsO.Range("B4").Style = "Hyperlink"
$wsO.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Write the new row3 (7188329f...) data
$wsO.Range("A3").Value2 = $newMd
$wsO.Range("B3").Value2 = "e2e\" + $newMd
$wsO.Range("B3").Style = "Hyperlink"
$wsO.Range("C3").Value2 = ".md"
$wsO.Range("E3").Value2 = "Handed back: in sync with en-US"
$wsO.Range("F3").Value2 = "Handed back: in sync with en-US"
$wsO.Range("G3").Value2 = "2016-08-30 22:48:24"
$wsO.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild hyperlinks for this sheet from scratch (existing hyperlink
# objects loaded from the source file can't be reliably edited/removed
# individually, so clear + re-add in the correct final layout).
$wsO.Cells.Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b718180687f249e0890af42b123bbd7eada1cd24/e2e/2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md", "", "", "e2e\2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B3"), $newMdUrlSrc, "", "", "e2e\" + $newMd) | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d270f66ff93f1dc1bffda9520d1b954c89c591f/e2e/b964db4f-3acf-4ace-850b-4e8e0e32472f.md", "", "", "e2e\b964db4f-3acf-4ace-850b-4e8e0e32472f.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn": columns A:P
# ---------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item(1)
$loZ.ListRows.Add() | Out-Null   # appends a blank row4, table range -> A1:P4

foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")) {
    $wsZ.Range($col + "4").Value2 = $wsZ.Range($col + "3").Value2
}
$wsZ.Range("A4").Style = "Hyperlink"
$wsZ.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("I4").Style = "Hyperlink"
$wsZ.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZ.Range("A3").Value2 = $newMd
$wsZ.Range("B3").Value2 = ".md"
$wsZ.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsZ.Range("D3").Value2 = "e2e"
$wsZ.Range("E3").Value2 = "ht"
$wsZ.Range("F3").Value2 = "True"
$wsZ.Range("G3").Value2 = $newZhXlf
$wsZ.Range("H3").Value2 = "2016-08-30 22:48:19"
$wsZ.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("I3").Value2 = $newMd
$wsZ.Range("I3").Style = "Hyperlink"
$wsZ.Range("J3").Value2 = $newZhXlf
$wsZ.Range("K3").Value2 = "2016-08-30 22:48:45"
$wsZ.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("L3").Value2 = ""
$wsZ.Range("M3").Value2 = "True"
$wsZ.Range("N3").Value2 = ""
$wsZ.Range("O3").Value2 = "False"
$wsZ.Range("P3").Value2 = ""
$wsZ.Range("A3").Style = "Hyperlink"

$wsZ.Cells.Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b718180687f249e0890af42b123bbd7eada1cd24/e2e/2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md", "", "", "2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1b0cb92a584e202539fda275b430648365b51d6b/e2e/2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md", "", "", "2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), $newMdUrlSrc, "", "", $newMd) | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I3"), $newMdUrlZh, "", "", $newMd) | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d270f66ff93f1dc1bffda9520d1b954c89c591f/e2e/b964db4f-3acf-4ace-850b-4e8e0e32472f.md", "", "", "b964db4f-3acf-4ace-850b-4e8e0e32472f.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3a5a703010ffff16198caeed267f3cfa062bf70e/e2e/b964db4f-3acf-4ace-850b-4e8e0e32472f.md", "", "", "b964db4f-3acf-4ace-850b-4e8e0e32472f.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de": columns A:P
# ---------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item(1)
$loD.ListRows.Add() | Out-Null   # appends a blank row4, table range -> A1:P4

foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")) {
    $wsD.Range($col + "4").Value2 = $wsD.Range($col + "3").Value2
}
$wsD.Range("A4").Style = "Hyperlink"
$wsD.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("I4").Style = "Hyperlink"
$wsD.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsD.Range("A3").Value2 = $newMd
$wsD.Range("B3").Value2 = ".md"
$wsD.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsD.Range("D3").Value2 = "e2e"
$wsD.Range("E3").Value2 = "ht"
$wsD.Range("F3").Value2 = "True"
$wsD.Range("G3").Value2 = $newDeXlf
$wsD.Range("H3").Value2 = "2016-08-30 22:48:24"
$wsD.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("I3").Value2 = $newMd
$wsD.Range("I3").Style = "Hyperlink"
$wsD.Range("J3").Value2 = $newDeXlf
$wsD.Range("K3").Value2 = "2016-08-30 22:48:53"
$wsD.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("L3").Value2 = ""
$wsD.Range("M3").Value2 = "True"
$wsD.Range("N3").Value2 = ""
$wsD.Range("O3").Value2 = "False"
$wsD.Range("P3").Value2 = ""
$wsD.Range("A3").Style = "Hyperlink"

$wsD.Cells.Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b718180687f249e0890af42b123bbd7eada1cd24/e2e/2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md", "", "", "2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7917d08037e7a1064198feee34350ab3766859d6/e2e/2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md", "", "", "2118c7ff-a5db-4c2b-bb59-20dd31a2f61a.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A3"), $newMdUrlSrc, "", "", $newMd) | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I3"), $newMdUrlDe, "", "", $newMd) | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d270f66ff93f1dc1bffda9520d1b954c89c591f/e2e/b964db4f-3acf-4ace-850b-4e8e0e32472f.md", "", "", "b964db4f-3acf-4ace-850b-4e8e0e32472f.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b89eddbf648aac2fd2d7ca8ec47680e82b52ce0d/e2e/b964db4f-3acf-4ace-850b-4e8e0e32472f.md", "", "", "b964db4f-3acf-4ace-850b-4e8e0e32472f.md") | Out-Null

Write-Output "Done."
